# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect the latest generated output (gh-pages build 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3264
$ws1.Range("F5").Value = 6902
$ws1.Range("F6").Value = 2166
$ws1.Range("F7").Value = 33
$ws1.Range("F8").Value = 82
$ws1.Range("F12").Value = 27
$ws1.Range("F13").Value = 155
$ws1.Range("F14").Value = 195

# Sheet "全部类型": row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3264
$ws4.Range("F6").Value = 6902
$ws4.Range("F7").Value = 2166
$ws4.Range("F8").Value = 33
$ws4.Range("F9").Value = 82
$ws4.Range("F13").Value = 27
$ws4.Range("F14").Value = 155
$ws4.Range("F15").Value = 195
